$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4, column F: PENDIENTE -> CORREGIDO ---
# (style moves from the "yellow / pending" look to the "green / corrected" look
# that's already used by several other rows in the sheet, e.g. F7/F8/F13-F16)
$ws.Range("F7").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F4").Value = "CORREGIDO"

# --- Row 18: fill in the new bug entry (previously an empty formatting-only row) ---

# Values are written in this order so new shared strings get appended to
# sharedStrings.xml in the same sequence as the target workbook.
$ws.Range("C18").Value = "Método obtenerUltimoTorneoDelUsurio en DAOTorneo: CAMBIAR!!!!"
$ws.Range("B18").Value = "Último Torneo del Usuario"
$ws.Range("E18").Value = "login.aspx"
$ws.Range("D18").Value = "Facu"
$ws.Range("A18").Value = 16
$ws.Range("F18").Value = "PENDIENTE"

# C18 gains word-wrap (same font/fill/border as before, just wrapped text).
$ws.Range("C18").WrapText = $true

# F18 takes on the same look as the other "PENDIENTE" cells (yellow fill,
# bold font, centered), but horizontally centered only (no vertical centering).
$ws.Range("F17").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("F18").VerticalAlignment = -4107
$ws.Range("F18").HorizontalAlignment = -4108
$ws.Range("F18").Value = "PENDIENTE"

# Row grows to fit the now two-line wrapped description.
$ws.Rows.Item(18).RowHeight = 30

# --- View state: scroll back to the top and move the active selection to I18 ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("I18").Select()
